# Update loading_percent values on the active sheet (case with 380 kV)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.44616561449082
$ws.Range("C2").Value = 3.782240501743753
$ws.Range("D2").Value = 6.017340911633282
$ws.Range("E2").Value = 12.05824923226895
$ws.Range("G2").Value = 55.25490992937782
$ws.Range("H2").Value = 20.11274330457283
$ws.Range("K2").Value = 9.894232612828695
$ws.Range("M2").Value = 14.35263990878995
$ws.Range("B3").Value = 10.29180863289551
$ws.Range("C3").Value = 3.595188410949474
$ws.Range("D3").Value = 5.904898994789359
$ws.Range("E3").Value = 11.80668895614482
$ws.Range("G3").Value = 54.22584642028062
$ws.Range("H3").Value = 19.9806304389955
$ws.Range("K3").Value = 9.820333475667416
$ws.Range("M3").Value = 14.23655979040004
$ws.Range("B4").Value = 10.20030792529846
$ws.Range("C4").Value = 3.474234378597765
$ws.Range("D4").Value = 5.836658637638327
$ws.Range("E4").Value = 11.65306579453098
$ws.Range("G4").Value = 53.59138883576521
$ws.Range("H4").Value = 19.90128047695091
$ws.Range("K4").Value = 9.77871330759433
$ws.Range("M4").Value = 14.16939250023881
$ws.Range("B5").Value = 10.16390179478366
$ws.Range("C5").Value = 3.4234272829312
$ws.Range("D5").Value = 5.809092178733639
$ws.Range("E5").Value = 11.59076375591952
$ws.Range("G5").Value = 53.33247363149819
$ws.Range("H5").Value = 19.86940645480706
$ws.Range("K5").Value = 9.762716629284329
$ws.Range("M5").Value = 14.14308076871943
$ws.Range("B6").Value = 10.15791144717146
$ws.Range("C6").Value = 3.41489974221901
$ws.Range("D6").Value = 5.804530606262926
$ws.Range("E6").Value = 11.58043937074157
$ws.Range("G6").Value = 53.28946718189618
$ws.Range("H6").Value = 19.86414217078619
$ws.Range("K6").Value = 9.760119148127407
$ws.Range("M6").Value = 14.13877645884574
$ws.Range("B7").Value = 10.19981329716413
$ws.Range("C7").Value = 3.473555294250494
$ws.Range("D7").Value = 5.836285832503584
$ws.Range("E7").Value = 11.65222422984822
$ws.Range("G7").Value = 53.58789812431935
$ws.Range("H7").Value = 19.90084872124389
$ws.Range("K7").Value = 9.7784936434439
$ws.Range("M7").Value = 14.16903332733062
$ws.Range("B8").Value = 10.39229871730925
$ws.Range("C8").Value = 3.719020966642048
$ws.Range("D8").Value = 5.978429594457465
$ws.Range("E8").Value = 11.97139166269467
$ws.Range("G8").Value = 54.90080100046625
$ws.Range("H8").Value = 20.06683441494997
$ws.Range("K8").Value = 9.867985549205542
$ws.Range("M8").Value = 14.31178000967453
$ws.Range("B9").Value = 10.79313394208279
$ws.Range("C9").Value = 4.151367216848188
$ws.Range("D9").Value = 6.261726394462505
$ws.Range("E9").Value = 12.60001383764587
$ws.Range("G9").Value = 57.44238446168934
$ws.Range("H9").Value = 20.40557179026659
$ws.Range("K9").Value = 10.07229520403238
$ws.Range("M9").Value = 14.62300389107792
$ws.Range("B10").Value = 11.09830804393692
$ws.Range("C10").Value = 4.43856154765359
$ws.Range("D10").Value = 6.470290787493203
$ws.Range("E10").Value = 13.05848481443277
$ws.Range("G10").Value = 59.27363385115727
$ws.Range("H10").Value = 20.6614176874745
$ws.Range("K10").Value = 10.23854001784828
$ws.Range("M10").Value = 14.86889166315867
$ws.Range("B11").Value = 11.23870052558101
$ws.Range("C11").Value = 4.562528780822918
$ws.Range("D11").Value = 6.564789225074999
$ws.Range("E11").Value = 13.26531929404437
$ws.Range("G11").Value = 60.09571504569831
$ws.Range("H11").Value = 20.77907771442029
$ws.Range("K11").Value = 10.31733820307569
$ws.Range("M11").Value = 14.98406953259473
$ws.Range("B12").Value = 11.29202680788726
$ws.Range("C12").Value = 4.608507604934637
$ws.Range("D12").Value = 6.600481347357827
$ws.Range("E12").Value = 13.34331585813736
$ws.Range("G12").Value = 60.40519019156362
$ws.Range("H12").Value = 20.82379419140007
$ws.Range("K12").Value = 10.34760351992841
$ws.Range("M12").Value = 15.02812491395783
$ws.Range("B13").Value = 11.28053578393008
$ws.Range("C13").Value = 4.598648210680108
$ws.Range("D13").Value = 6.592799103664258
$ws.Range("E13").Value = 13.32653366100309
$ws.Range("G13").Value = 60.33862441509054
$ws.Range("H13").Value = 20.8141568848335
$ws.Range("K13").Value = 10.3410668761122
$ws.Range("M13").Value = 15.01861783729153
$ws.Range("B14").Value = 11.24308480364242
$ws.Range("C14").Value = 4.566330844092743
$ws.Range("D14").Value = 6.567727699808787
$ws.Range("E14").Value = 13.27174311091498
$ws.Range("G14").Value = 60.12121381113609
$ws.Range("H14").Value = 20.78275346338094
$ws.Range("K14").Value = 10.31981978315734
$ws.Range("M14").Value = 14.98768539522145
$ws.Range("B15").Value = 11.22016435115473
$ws.Range("C15").Value = 4.546409762091822
$ws.Range("D15").Value = 6.552357614810131
$ws.Range("E15").Value = 13.23813744632422
$ws.Range("G15").Value = 59.98779805633724
$ws.Range("H15").Value = 20.76353828969275
$ws.Range("K15").Value = 10.30685991715238
$ws.Range("M15").Value = 14.96879455332388
$ws.Range("B16").Value = 11.08915921080572
$ws.Range("C16").Value = 4.430325184303578
$ws.Range("D16").Value = 6.464104189467407
$ws.Range("E16").Value = 13.04492610678846
$ws.Range("G16").Value = 59.21966640348052
$ws.Range("H16").Value = 20.65375209351249
$ws.Range("K16").Value = 10.23345166134121
$ws.Range("M16").Value = 14.86142816455099
$ws.Range("B17").Value = 11.00914816431776
$ws.Range("C17").Value = 4.357396755556538
$ws.Range("D17").Value = 6.409839044840251
$ws.Range("E17").Value = 12.92589747612447
$ws.Range("G17").Value = 58.74545062527447
$ws.Range("H17").Value = 20.58671263983488
$ws.Range("K17").Value = 10.18920953902871
$ws.Range("M17").Value = 14.79638669590753
$ws.Range("B18").Value = 10.96327911589576
$ws.Range("C18").Value = 4.314821760104699
$ws.Range("D18").Value = 6.378593699177625
$ws.Range("E18").Value = 12.85727731343768
$ws.Range("G18").Value = 58.47167663127164
$ws.Range("H18").Value = 20.54827462395438
$ws.Range("K18").Value = 10.16406312073272
$ws.Range("M18").Value = 14.75929158855059
$ws.Range("B19").Value = 10.94777654467223
$ws.Range("C19").Value = 4.300298815503088
$ws.Range("D19").Value = 6.368009999434339
$ws.Range("E19").Value = 12.83401904078981
$ws.Range("G19").Value = 58.37881459696979
$ws.Range("H19").Value = 20.53528167324315
$ws.Range("K19").Value = 10.15560148514055
$ws.Range("M19").Value = 14.74678706491566
$ws.Range("B20").Value = 11.01765028814169
$ws.Range("C20").Value = 4.365225194896342
$ws.Range("D20").Value = 6.415619406770199
$ws.Range("E20").Value = 12.93858522003759
$ws.Range("G20").Value = 58.79603889645099
$ws.Range("H20").Value = 20.59383670854233
$ws.Range("K20").Value = 10.19388830509117
$ws.Range("M20").Value = 14.80327813827161
$ws.Range("B21").Value = 11.25408114052739
$ws.Range("C21").Value = 4.575849452570534
$ws.Range("D21").Value = 6.575094579469251
$ws.Range("E21").Value = 13.28784590898699
$ws.Range("G21").Value = 60.18512416600637
$ws.Range("H21").Value = 20.79197321287999
$ws.Range("K21").Value = 10.326049247205
$ws.Range("M21").Value = 14.99675936533563
$ws.Range("B22").Value = 11.40952281742794
$ws.Range("C22").Value = 4.707879489460334
$ws.Range("D22").Value = 6.678766061987069
$ws.Range("E22").Value = 13.51416767637662
$ws.Range("G22").Value = 61.08220226495835
$ws.Range("H22").Value = 20.9223958738272
$ws.Range("K22").Value = 10.41489356264529
$ws.Range("M22").Value = 15.12575737595854
$ws.Range("B23").Value = 11.32649672733564
$ws.Range("C23").Value = 4.637928363668578
$ws.Range("D23").Value = 6.623497455506998
$ws.Range("E23").Value = 13.39357790246875
$ws.Range("G23").Value = 60.60448020331387
$ws.Range("H23").Value = 20.85270922271247
$ws.Range("K23").Value = 10.36725984501127
$ws.Range("M23").Value = 15.05668830561569
$ws.Range("B24").Value = 11.01380606708896
$ws.Range("C24").Value = 4.361687973892409
$ws.Range("D24").Value = 6.413006250497877
$ws.Range("E24").Value = 12.93284967376543
$ws.Range("G24").Value = 58.77317147375584
$ws.Range("H24").Value = 20.59061559182435
$ws.Range("K24").Value = 10.19177213257436
$ws.Range("M24").Value = 14.800161585789
$ws.Range("B25").Value = 10.68256247577874
$ws.Range("C25").Value = 4.039721910978551
$ws.Range("D25").Value = 6.184846114730412
$ws.Range("E25").Value = 12.43019486242222
$ws.Range("G25").Value = 56.76003424880606
$ws.Range("H25").Value = 20.31262747223261
$ws.Range("K25").Value = 10.014080264389
$ws.Range("M25").Value = 14.5356433801687
